$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with the latest scraped figures. The price column is formatted as Text
# first so Excel keeps the exact original string (leading/trailing zeros,
# "thousand dot" grouping, etc.) instead of silently coercing it to a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.442.63"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.02"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.46"
$ws.Range("E5").Value = "  -7.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5243"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3239"
$ws.Range("E8").Value = "  -8.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06746"
$ws.Range("E9").Value = "  -4.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.89"
$ws.Range("E10").Value = "  -7.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7725"
$ws.Range("E11").Value = "  -5.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07680"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.852.45"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.05"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.022"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.13"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007874"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.466.27"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.078.38"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.530"
$ws.Range("E22").Value = "  -5.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.437"
$ws.Range("E23").Value = "  -7.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.916"
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.03"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.93"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.34"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.182"
$ws.Range("E30").Value = "  -4.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08768"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.101"
$ws.Range("E32").Value = "  -6.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04839"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.133"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.856"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6878"
$ws.Range("E36").Value = "  -7.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.101"
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01791"
$ws.Range("E38").Value = "  -4.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.215"
$ws.Range("E39").Value = "  -8.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4933"
$ws.Range("E40").Value = "  -7.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.88"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8985"
$ws.Range("E42").Value = "  -8.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.180"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.748"
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4195"
$ws.Range("E46").Value = "  -8.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1260"
$ws.Range("E47").Value = "  -7.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.048"
$ws.Range("E48").Value = "  -4.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05875"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.39"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.24"
$ws.Range("E51").Value = "  -4.25%  "
